$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source file has an empty placeholder string cell at F1 (t="s" with no
# value). Re-assert it as empty so the COM round-trip doesn't resolve it to
# shared string index 0 ("number").
$ws.Range("F1").Value = ""

# Populate the newly added "duplicate_image_filename" values (column E) for
# the data rows, matching the "NA" placeholder added by the commit.
foreach ($r in 2..21) {
    $ws.Range("E$r").Value = "NA"
}
